$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1573.4688
$ws.Range("I138").Value = 1075.0454
$ws.Range("J138").Value = 2670
$ws.Range("K138").Value = 3225.1362
$ws.Range("L138").Value = 8010
$ws.Range("M138").Value = 1914.8638
$ws.Range("N138").Value = -18290

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 6822
$ws.Range("I41").Value = 3269.75
$ws.Range("J41").Value = 21031
$ws.Range("K41").Value = 3269.75
$ws.Range("L41").Value = 21031
$ws.Range("M41").Value = -2855.75
$ws.Range("N41").Value = -21859
$ws.Range("H132").Value = 8046.684
$ws.Range("I132").Value = 12151.223
$ws.Range("K132").Value = 36453.669
$ws.Range("M132").Value = -33923.669

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 6065.4443
$ws.Range("I36").Value = 1709.6
$ws.Range("J36").Value = 11510.25
$ws.Range("K36").Value = 1709.6
$ws.Range("L36").Value = 11510.25
$ws.Range("M36").Value = -1175.6
$ws.Range("N36").Value = -12578.25
$ws.Range("H49").Value = 10000
$ws.Range("J49").Value = 10000
$ws.Range("L49").Value = 10000
$ws.Range("N49").Value = -10478
$ws.Range("H51").Value = 43437.5
$ws.Range("I51").Value = 20000
$ws.Range("J51").Value = 66875
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 66875
$ws.Range("M51").Value = -19509
$ws.Range("N51").Value = -67857
$ws.Range("H80").Value = 332.0909
$ws.Range("J80").Value = 170.94118
$ws.Range("L80").Value = 170.94118
$ws.Range("N80").Value = -2166.94118
$ws.Range("H83").Value = 332.0909
$ws.Range("J83").Value = 170.94118
$ws.Range("L83").Value = 854.7059
$ws.Range("N83").Value = -10838.7059
$ws.Range("H107").Value = 1505.95
$ws.Range("I107").Value = 1018.7273
$ws.Range("J107").Value = 2101.4443
$ws.Range("K107").Value = 1018.7273
$ws.Range("L107").Value = 2101.4443
$ws.Range("M107").Value = 901.2727
$ws.Range("N107").Value = -5941.4443
$ws.Range("H125").Value = 55566.668
$ws.Range("J125").Value = 55566.668
$ws.Range("L125").Value = 55566.668
$ws.Range("N125").Value = -65406.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 43296.168
$ws.Range("J20").Value = 43296.168
$ws.Range("L20").Value = 43296.168
$ws.Range("N20").Value = -43768.168
$ws.Range("H30").Value = 43296.168
$ws.Range("J30").Value = 43296.168
$ws.Range("L30").Value = 43296.168
$ws.Range("N30").Value = -43478.168
$ws.Range("H58").Value = 1033.0769
$ws.Range("I58").Value = 1072.25
$ws.Range("K58").Value = 1072.25
$ws.Range("M58").Value = -869.25
$ws.Range("H99").Value = 2349
$ws.Range("I99").Value = 2148.6667
$ws.Range("J99").Value = 2499.25
$ws.Range("K99").Value = 2148.6667
$ws.Range("L99").Value = 2499.25
$ws.Range("M99").Value = -650.6667000000002
$ws.Range("N99").Value = -5495.25
$ws.Range("H126").Value = 2349
$ws.Range("I126").Value = 2148.6667
$ws.Range("J126").Value = 2499.25
$ws.Range("K126").Value = 6446.000100000001
$ws.Range("L126").Value = 7497.75
$ws.Range("M126").Value = -3976.000100000001
$ws.Range("N126").Value = -12437.75
$ws.Range("H128").Value = 43296.168
$ws.Range("J128").Value = 43296.168
$ws.Range("L128").Value = 43296.168
$ws.Range("N128").Value = -53256.168
$ws.Range("H129").Value = 39037.125
$ws.Range("I129").Value = 10000
$ws.Range("K129").Value = 10000
$ws.Range("M129").Value = -5000
$ws.Range("H132").Value = 2282.4092
$ws.Range("I132").Value = 1453.8235
$ws.Range("J132").Value = 5099.6
$ws.Range("K132").Value = 4361.470499999999
$ws.Range("L132").Value = 15298.8
$ws.Range("M132").Value = -1831.470499999999
$ws.Range("N132").Value = -20358.8
$ws.Range("H133").Value = 35363.637
$ws.Range("J133").Value = 35363.637
$ws.Range("L133").Value = 35363.637
$ws.Range("N133").Value = -40423.637
$ws.Range("H136").Value = 1033.0769
$ws.Range("I136").Value = 1072.25
$ws.Range("K136").Value = 3216.75
$ws.Range("M136").Value = -666.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1758.909
$ws.Range("I34").Value = 749.5
$ws.Range("J34").Value = 1983.2222
$ws.Range("K34").Value = 2248.5
$ws.Range("L34").Value = 5949.6666
$ws.Range("M34").Value = -2164.5
$ws.Range("N34").Value = -6117.6666
$ws.Range("H39").Value = 2412.375
$ws.Range("J39").Value = 2412.375
$ws.Range("L39").Value = 7237.125
$ws.Range("N39").Value = -7825.125
$ws.Range("H55").Value = 3484.25
$ws.Range("J55").Value = 3981.1
$ws.Range("L55").Value = 11943.3
$ws.Range("N55").Value = -12297.3
$ws.Range("H63").Value = 3996.6667
$ws.Range("I63").Value = 3996.6667
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 11990.0001
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -11241.0001
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 3996.6667
$ws.Range("I66").Value = 3996.6667
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 35970.0003
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -32226.0003
$ws.Range("N66").ClearContents()
$ws.Range("H68").Value = 1133.025
$ws.Range("I68").Value = 1075.4814
$ws.Range("J68").Value = 1252.5385
$ws.Range("K68").Value = 3226.4442
$ws.Range("L68").Value = 3757.6155
$ws.Range("M68").Value = -2415.4442
$ws.Range("N68").Value = -5379.6155
$ws.Range("H71").Value = 1133.025
$ws.Range("I71").Value = 1075.4814
$ws.Range("J71").Value = 1252.5385
$ws.Range("K71").Value = 9679.3326
$ws.Range("L71").Value = 11272.8465
$ws.Range("M71").Value = -5623.3326
$ws.Range("N71").Value = -19384.8465
$ws.Range("H107").Value = 647.8333
$ws.Range("J107").Value = 691.3125
$ws.Range("L107").Value = 2073.9375
$ws.Range("N107").Value = -5913.9375
$ws.Range("H116").Value = 11191.728
$ws.Range("I116").Value = 15415.571
$ws.Range("J116").Value = 3800
$ws.Range("K116").Value = 46246.713
$ws.Range("L116").Value = 11400
$ws.Range("M116").Value = -42804.713
$ws.Range("N116").Value = -18284
$ws.Range("H122").Value = 849.8461
$ws.Range("I122").Value = 423.70587
$ws.Range("J122").Value = 1654.7778
$ws.Range("K122").Value = 3813.35283
$ws.Range("L122").Value = 14893.0002
$ws.Range("M122").Value = -1363.35283
$ws.Range("N122").Value = -19793.0002
$ws.Range("H129").Value = 1284.2
$ws.Range("I129").Value = 541.6667
$ws.Range("J129").Value = 1779.2222
$ws.Range("K129").Value = 1625.0001
$ws.Range("L129").Value = 5337.6666
$ws.Range("M129").Value = 3374.9999
$ws.Range("N129").Value = -15337.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2454.5625
$ws.Range("I7").Value = 1996.75
$ws.Range("J7").Value = 2912.375
$ws.Range("K7").Value = 1996.75
$ws.Range("L7").Value = 2912.375
$ws.Range("M7").Value = -1884.75
$ws.Range("N7").Value = -3136.375
$ws.Range("H125").Value = 66715
$ws.Range("J125").Value = 66715
$ws.Range("L125").Value = 66715
$ws.Range("N125").Value = -76555
$ws.Range("H126").Value = 2454.5625
$ws.Range("I126").Value = 1996.75
$ws.Range("J126").Value = 2912.375
$ws.Range("K126").Value = 5990.25
$ws.Range("L126").Value = 8737.125
$ws.Range("M126").Value = -3520.25
$ws.Range("N126").Value = -13677.125
$ws.Range("H132").Value = 7674.029
$ws.Range("I132").Value = 12383.444
$ws.Range("J132").Value = 2687.5881
$ws.Range("K132").Value = 37150.33199999999
$ws.Range("L132").Value = 8062.7643
$ws.Range("M132").Value = -34620.33199999999
$ws.Range("N132").Value = -13122.7643

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 8000
$ws.Range("I51").Value = 8000
$ws.Range("K51").Value = 8000
$ws.Range("M51").Value = -7490
$ws.Range("H126").Value = 926.6923
$ws.Range("I126").Value = 971.6957
$ws.Range("J126").Value = 581.6667
$ws.Range("K126").Value = 2915.0871
$ws.Range("L126").Value = 1745.0001
$ws.Range("M126").Value = -445.0870999999997
$ws.Range("N126").Value = -6685.0001
$ws.Range("H132").Value = 1268.0834
$ws.Range("I132").Value = 880.05554
$ws.Range("J132").Value = 2432.1667
$ws.Range("K132").Value = 2640.16662
$ws.Range("L132").Value = 7296.500100000001
$ws.Range("M132").Value = -110.16662
$ws.Range("N132").Value = -12356.5001
